# Updated symbol list on Mon Jan  2 03:45:30 UTC 2023 with GitHub Actions
#
# This script refreshes the "Price" (column D) and "Volume(1h)" (column E)
# figures on the active worksheet for the coin rows whose quotes changed.
# The values are written as literal text (matching the workbook's original
# inline-string cell type) rather than as numbers/percentages, so each cell
# is temporarily switched to the "Text" number format while the value is
# assigned and then restored to the sheet's normal style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "242.97" },
    @{ Cell = "E2"; Value = "-0.95%" },
    @{ Cell = "D3"; Value = "29.66" },
    @{ Cell = "E3"; Value = "11.79%" },
    @{ Cell = "D4"; Value = "5.132" },
    @{ Cell = "E4"; Value = "0.52%" },
    @{ Cell = "D5"; Value = "0.05655" },
    @{ Cell = "E5"; Value = "1.20%" },
    @{ Cell = "D6"; Value = "6.499" },
    @{ Cell = "E6"; Value = "0.38%" },
    @{ Cell = "D7"; Value = "0.8274" },
    @{ Cell = "E7"; Value = "1.24%" },
    @{ Cell = "D8"; Value = "0.8619" },
    @{ Cell = "E8"; Value = "2.85%" },
    @{ Cell = "D9"; Value = "0.1329" },
    @{ Cell = "E9"; Value = "0.08%" },
    @{ Cell = "D10"; Value = "0.06913" },
    @{ Cell = "E10"; Value = "-1.04%" },
    @{ Cell = "D11"; Value = "0.02855" },
    @{ Cell = "E11"; Value = "-1.07%" },
    @{ Cell = "E12"; Value = "0.06%" },
    @{ Cell = "D13"; Value = "0.001523" },
    @{ Cell = "E13"; Value = "1.06%" },
    @{ Cell = "D14"; Value = "0.04149" },
    @{ Cell = "E14"; Value = "-9.42%" },
    @{ Cell = "D15"; Value = "0.0005976" },
    @{ Cell = "E15"; Value = "0.28%" },
    @{ Cell = "D16"; Value = "0.006213" },
    @{ Cell = "E16"; Value = "1.10%" },
    @{ Cell = "E17"; Value = "-3.12%" },
    @{ Cell = "D18"; Value = "3.011" },
    @{ Cell = "E18"; Value = "-0.84%" },
    @{ Cell = "D19"; Value = "2.220" },
    @{ Cell = "E19"; Value = "1.69%" },
    @{ Cell = "D21"; Value = "0.03240" },
    @{ Cell = "E21"; Value = "4.42%" },
    @{ Cell = "D22"; Value = "0.1294" },
    @{ Cell = "E22"; Value = "-0.38%" },
    @{ Cell = "D23"; Value = "3.613" },
    @{ Cell = "E23"; Value = "-3.68%" },
    @{ Cell = "D24"; Value = "0.1373" },
    @{ Cell = "E24"; Value = "-0.05%" },
    @{ Cell = "D25"; Value = "0.001209" },
    @{ Cell = "E25"; Value = "-2.78%" },
    @{ Cell = "D26"; Value = "0.004446" },
    @{ Cell = "E26"; Value = "-1.51%" },
    @{ Cell = "E27"; Value = "22.84%" },
    @{ Cell = "D28"; Value = "0.0001403" },
    @{ Cell = "E28"; Value = "0.51%" },
    @{ Cell = "E40"; Value = "1.77%" },
    @{ Cell = "D41"; Value = "0.005776" },
    @{ Cell = "E41"; Value = "67.30%" },
    @{ Cell = "D42"; Value = "0.1054" },
    @{ Cell = "E42"; Value = "-23.13%" },
    @{ Cell = "E43"; Value = "-12.16%" },
    @{ Cell = "D45"; Value = "0.00005096" },
    @{ Cell = "E45"; Value = "-4.54%" },
    @{ Cell = "D46"; Value = "0.00000000749" },
    @{ Cell = "E46"; Value = "-0.06%" },
    @{ Cell = "D47"; Value = "0.1009" },
    @{ Cell = "E47"; Value = "-7.39%" },
    @{ Cell = "D48"; Value = "0.002741" },
    @{ Cell = "E48"; Value = "7.41%" },
    @{ Cell = "E49"; Value = "-0.06%" },
    @{ Cell = "E50"; Value = "-0.06%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
